$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    # Force the cell to store a literal text value (never an auto-converted
    # number/date), then strip the quote-prefix formatting the leading
    # apostrophe trick leaves behind so the cell is left unstyled again.
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).ClearFormats()
}

# Row 2 - Bitcoin
Set-TextCell "D2" "66.667.20"
Set-TextCell "E2" "  -4.41%  "

# Row 3 - Ethereum
Set-TextCell "D3" "3.321.23"
Set-TextCell "E3" "  -1.81%  "

# Row 4 - TetherUSD
Set-TextCell "E4" "  +0.02%  "

# Row 5 - BNB
Set-TextCell "D5" "572.98"
Set-TextCell "E5" "  -3.61%  "

# Row 6 - Solana
Set-TextCell "D6" "182.24"
Set-TextCell "E6" "  -5.71%  "

# Row 7 - USDC
Set-TextCell "E7" "  +0.03%  "

# Row 8 - XRP
Set-TextCell "E8" "  -1.78%  "

# Row 9 - Dogecoin
Set-TextCell "E9" "  -4.33%  "

# Row 10 - Toncoin
Set-TextCell "E10" "  -1.97%  "

# Row 11 - Cardano
Set-TextCell "D11" "0.403"
Set-TextCell "E11" "  -4.90%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextCell "D12" "3.901.13"
Set-TextCell "E12" "  -1.75%  "

# Row 13 - TRON
Set-TextCell "E13" "  -0.93%  "

# Row 14 - Avalanche
Set-TextCell "D14" "27.16"
Set-TextCell "E14" "  -5.64%  "

# Row 15 - WrappedBTC
Set-TextCell "D15" "66.746.66"
Set-TextCell "E15" "  -4.29%  "

# Row 16 - ShibaInu
Set-TextCell "E16" "  -3.30%  "

# Row 17 - WrappedEther
Set-TextCell "D17" "3.297.08"
Set-TextCell "E17" "  -2.48%  "

# Row 18 - Chainlink
Set-TextCell "D18" "13.64"
Set-TextCell "E18" "  -1.17%  "

# Row 19 - was BitcoinCash, now Polkadot
Set-TextCell "B19" "Polkadot"
Set-TextCell "C19" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell "D19" "5.69"
Set-TextCell "E19" "  -2.88%  "

# Row 20 - was Polkadot, now BitcoinCash
Set-TextCell "B20" "BitcoinCash"
Set-TextCell "C20" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextCell "D20" "432.04"
Set-TextCell "E20" "  -4.27%  "

# Row 21 - Uniswap
Set-TextCell "E21" "  -2.78%  "

# Row 22 - Litecoin
Set-TextCell "D22" "73.53"
Set-TextCell "E22" "  -0.02%  "

# Row 23 - Dai
Set-TextCell "E23" "  -0.19%  "

# Row 24 - Polygon
Set-TextCell "D24" "0.521"
Set-TextCell "E24" "  +0.32%  "

# Row 25 - PEPE
Set-TextCell "E25" "  -3.25%  "

# Row 26 - Kaspa
Set-TextCell "E26" "  -1.21%  "

# Row 27 - InternetComputer(DFINITY)
Set-TextCell "D27" "9.10"
Set-TextCell "E27" "  -5.54%  "

# Row 28 - Binance-PegBSC-USD
Set-TextCell "E28" "  -2.37%  "

# Row 29 - PancakeSwap
Set-TextCell "E29" "  -2.38%  "

# Row 30 - EthereumClassic
Set-TextCell "D30" "22.81"
Set-TextCell "E30" "  -2.02%  "

# Row 31 - NEARProtocol
Set-TextCell "D31" "5.32"
Set-TextCell "E31" "  -5.92%  "

# Row 32 - USDe
Set-TextCell "E32" "  +0.01%  "

# Row 33 - was Fetch.AI, now Aptos
Set-TextCell "B33" "Aptos"
Set-TextCell "C33" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextCell "D33" "6.84"
Set-TextCell "E33" "  -3.10%  "

# Row 34 - was Aptos, now Fetch.AI
Set-TextCell "B34" "Fetch.AI"
Set-TextCell "C34" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextCell "D34" "1.24"
Set-TextCell "E34" "  -4.47%  "

# Row 35 - ImmutableX
Set-TextCell "E35" "  -1.44%  "

# Row 36 - Monero
Set-TextCell "D36" "159.85"
Set-TextCell "E36" "  -3.02%  "

# Row 37 - EnergySwap
Set-TextCell "D37" "27.30"
Set-TextCell "E37" "  -0.04%  "

# Row 38 - Stacks
Set-TextCell "E38" "  -4.86%  "

# Row 39 - Maker
Set-TextCell "D39" "2.835.54"
Set-TextCell "E39" "  +3.16%  "

# Row 40 - Mantle
Set-TextCell "D40" "0.790"
Set-TextCell "E40" "  -4.16%  "

# Row 41 - Filecoin
Set-TextCell "D41" "4.44"
Set-TextCell "E41" "  -3.83%  "

# Row 42 - RenderToken
Set-TextCell "E42" "  -4.94%  "

# Row 43 - Hedera
Set-TextCell "E43" "  -2.28%  "

# Row 44 - OKB
Set-TextCell "D44" "40.14"
Set-TextCell "E44" "  -1.67%  "

# Row 45 - InjectiveProtocol
Set-TextCell "D45" "24.47"
Set-TextCell "E45" "  -4.74%  "

# Row 46 - dogwifhat
Set-TextCell "E46" "  -7.90%  "

# Row 47 - Bittensor
Set-TextCell "D47" "324.57"
Set-TextCell "E47" "  -6.17%  "

# Row 48 - VeChain
Set-TextCell "E48" "  -5.13%  "

# Row 49 - ONDO
Set-TextCell "D49" "0.988"
Set-TextCell "E49" "  -2.45%  "

# Row 50 - Cosmos
Set-TextCell "D50" "6.18"
Set-TextCell "E50" "  -2.81%  "

# Row 51 - Stellar
Set-TextCell "E51" "  -1.64%  "
